$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new "Item starts with M" entry ---
$ws.Range("A8").Value2 = "Item starts with M use:"
$ws.Range("B8").Value2 = "F-825-1052C Inspection form CMD PN M.C00.00C.xlsx"
$ws.Rows.Item(8).RowHeight = 28.8

# --- Rows 10 & 11: clear the old CofC rows that used to live here ---
$ws.Range("A10:B10").ClearContents()
$ws.Rows.Item(10).EntireRow.AutoFit()

$ws.Range("A11:B11").ClearContents()
$ws.Rows.Item(11).EntireRow.AutoFit()

# --- Rows 14 & 15: replace "Document Type/Where Stored" table with the CofC info ---
$ws.Range("A14").Value2 = "CofC for sterile "
$ws.Range("B14").Value2 = "F-825-015B S3D CofC Carlsmed Sterile Product.docx"
$ws.Range("A14").VerticalAlignment = -4160
$ws.Range("B14").VerticalAlignment = -4160

$ws.Range("A15").Value2 = "CofC for Non-Sterile"
$ws.Range("B15").Value2 = "F-825-008Q S3D CofC (Carlsmed).docx"
$ws.Range("A15").VerticalAlignment = -4160
$ws.Range("B15").VerticalAlignment = -4160

# --- Rows 16, 17, 18: clear old CofC/Final Inspection/Dim & Meas rows ---
$ws.Range("A16:B16").ClearContents()
$ws.Rows.Item(16).EntireRow.AutoFit()

$ws.Range("A17:B17").ClearContents()
$ws.Rows.Item(17).EntireRow.AutoFit()

$ws.Range("A18:B18").ClearContents()
$ws.Rows.Item(18).EntireRow.AutoFit()

# --- Column B a bit wider ---
$ws.Columns.Item(2).ColumnWidth = 32.8333333333333

# --- Selection moves up one row ---
$ws.Range("B11").Select()
